$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.948
$ws.Range("C4").Value = -12.78
$ws.Range("E4").Value = 12.938

$ws.Range("C5").Value = -12.413

$ws.Range("A7").Value = -20.987

$ws.Range("C8").Value = -12.34

$ws.Range("E9").Value = 13.021

$ws.Range("A16").Value = -21.355
$ws.Range("C16").Value = -12.346

$ws.Range("E18").Value = 13.268
